$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B5: the "RunTest" flag for CT 04 flips from Yes to No
$ws.Range("B5").Value = "No"

# H2 and H5: the "vOutData" timestamp moves from 30/05/2020 to 01/06/2020.
# A plain Range.Value assignment would make the host auto-detect the
# date-shaped text and store it as a date serial (changing both the cell
# type and its number-format style), so instead we write it as a text
# formula and then collapse the formula down to its literal value with a
# values-only paste. That keeps the cell a shared-string text cell and
# leaves its original style untouched, just like the source edit.
$ws.Range("H2").Formula = '="01/06/2020"'
$ws.Range("H2").Copy()
$ws.Range("H2").PasteSpecial(-4163)

$ws.Range("H5").Formula = '="01/06/2020"'
$ws.Range("H5").Copy()
$ws.Range("H5").PasteSpecial(-4163)

# Match the author's final selection/active-cell position.
$ws.Range("D11").Select()

$wb.Save()
